$d = $word.ActiveDocument

$pairs = @(
    ,@('2023-04-27 Thursday', '2023-04-28 Friday')
    ,@('64×32=', '16×100=')
    ,@('12×39=', '18×37=')
    ,@('67×83=', '44×62=')
    ,@('22×65=', '33×16=')
    ,@('64×27=', '42×27=')
    ,@('77×27=', '63×54=')
    ,@('93×72=', '27×99=')
    ,@('66×86=', '19×24=')
    ,@('26×72=', '55×45=')
    ,@('70×97=', '26×73=')
    ,@('29×62=', '13×65=')
    ,@('30×53=', '38×54=')
    ,@('56×28=', '84×92=')
    ,@('19×61=', '12×86=')
    ,@('59×28=', '42×41=')
    ,@('92×97=', '33×94=')
    ,@('54×70=', '20×20=')
    ,@('59×75=', '76×53=')
    ,@('38×64=', '61×10=')
    ,@('88×67=', '75×66=')
    ,@('85×15=', '66×59=')
    ,@('27×35=', '55×91=')
    ,@('43×15=', '88×85=')
    ,@('50×46=', '77×62=')
    ,@('20×18=', '41×47=')
    ,@('93×94=', '84×45=')
    ,@('59×71=', '87×99=')
    ,@('21×98=', '80×99=')
    ,@('60×99=', '45×37=')
    ,@('83×44=', '73×39=')
    ,@('38×97=', '97×82=')
    ,@('60×93=', '57×88=')
    ,@('91×84=', '41×95=')
    ,@('86×32=', '95×39=')
    ,@('47×74=', '20×78=')
    ,@('22×87=', '27×18=')
    ,@('51×35=', '40×41=')
    ,@('33×40=', '88×12=')
    ,@('62×25=', '80×56=')
    ,@('38×47=', '14×27=')
    ,@('77×100=', '67×92=')
    ,@('68×45=', '73×82=')
    ,@('11×31=', '65×43=')
    ,@('82×32=', '31×57=')
    ,@('83×45=', '21×31=')
    ,@('81×100=', '38×32=')
    ,@('61×12=', '70×22=')
    ,@('13×92=', '68×74=')
    ,@('26×33=', '16×82=')
    ,@('30×31=', '46×36=')
    ,@('12×16=', '61×14=')
    ,@('79×72=', '38×79=')
    ,@('74×89=', '25×61=')
    ,@('66×12=', '24×21=')
    ,@('11×65=', '68×14=')
    ,@('90×27=', '71×60=')
    ,@('42×74=', '48×27=')
    ,@('50×23=', '65×58=')
    ,@('92×64=', '62×22=')
    ,@('34×37=', '57×46=')
    ,@('68×10=', '40×100=')
    ,@('33×19=', '91×81=')
    ,@('16×93=', '95×48=')
    ,@('91×62=', '11×46=')
    ,@('22×77=', '38×11=')
    ,@('77×56=', '30×73=')
    ,@('47×20=', '59×24=')
    ,@('17×75=', '57×56=')
    ,@('77×97=', '29×77=')
    ,@('34×79=', '66×83=')
    ,@('47×68=', '38×60=')
    ,@('95×81=', '84×54=')
    ,@('69×16=', '31×43=')
    ,@('31×20=', '49×88=')
    ,@('90×48=', '70×15=')
    ,@('19×99=', '57×59=')
    ,@('69×66=', '51×83=')
    ,@('51×27=', '79×42=')
    ,@('17×32=', '68×100=')
    ,@('51×96=', '96×96=')
    ,@('44×83=', '20×17=')
    ,@('66×28=', '54×43=')
    ,@('97×70=', '56×29=')
    ,@('56×35=', '23×70=')
    ,@('23×62=', '42×16=')
    ,@('93×28=', '97×41=')
    ,@('97×23=', '26×16=')
    ,@('67×28=', '39×57=')
    ,@('13×60=', '93×22=')
    ,@('16×53=', '39×35=')
    ,@('12×79=', '31×89=')
    ,@('41×83=', '17×50=')
    ,@('39×24=', '98×89=')
    ,@('80×79=', '89×29=')
    ,@('51×55=', '85×92=')
    ,@('63×65=', '10×43=')
    ,@('17×62=', '70×89=')
    ,@('93×59=', '21×57=')
    ,@('92×22=', '51×39=')
    ,@('50×41=', '21×88=')
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
